$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.403.41"
$ws.Range("D3").Value = "1.548.21"
$ws.Range("E3").Value = "  -2.22%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'209.69"
$ws.Range("E5").Value = "  -1.85%  "
$ws.Range("E6").Value = "  -2.26%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'24.00"
$ws.Range("E8").Value = "  -0.50%  "
$ws.Range("E9").Value = "  -2.05%  "
$ws.Range("E10").Value = "  -1.80%  "
$ws.Range("D11").Value = "'0.0887"
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").Value = "1.770.07"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("D13").Value = "1.550.10"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("D14").Value = "28.371.99"
$ws.Range("E14").Value = "  -0.88%  "
$ws.Range("E15").Value = "  -2.78%  "
$ws.Range("D16").Value = "'0.508"
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").Value = "'60.80"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "'228.73"
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("D19").Value = "'7.32"
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("D20").Value = "0.0₃0670"
$ws.Range("E20").Value = "  -3.25%  "
$ws.Range("E21").Value = "  -0.10%  "
$ws.Range("D22").Value = "'3.88"
$ws.Range("E22").Value = "  -0.90%  "
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("D25").Value = "'151.21"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "'14.74"
$ws.Range("E26").Value = "  -2.18%  "
$ws.Range("D27").Value = "'0.102"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "'6.22"
$ws.Range("E29").Value = "  -3.77%  "
$ws.Range("D30").Value = "'0.0466"
$ws.Range("E30").Value = "  -3.37%  "
$ws.Range("E31").Value = "  -5.27%  "
$ws.Range("D33").Value = "1.381.07"
$ws.Range("E33").Value = "  -1.37%  "
$ws.Range("D34").Value = "'2.99"
$ws.Range("E34").Value = "  -3.94%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("E36").Value = "  -3.74%  "
$ws.Range("D37").Value = "'2.31"
$ws.Range("E37").Value = "  -2.23%  "
$ws.Range("D38").Value = "'2.56"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("E39").Value = "  -3.01%  "
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("B41").Value = "ImmutableX"
$ws.Range("C41").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D41").Value = "'0.508"
$ws.Range("E41").Value = "  -2.85%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "'0.767"
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("D44").Value = "'0.0456"
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("D45").Value = "'5.33"
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("D46").Value = "'61.93"
$ws.Range("D47").Value = "1.683.27"
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("D48").Value = "'0.873"
$ws.Range("E48").Value = "  -9.21%  "
$ws.Range("D49").Value = "'43.96"
$ws.Range("E49").Value = "  +9.54%  "
$ws.Range("D50").Value = "'85.45"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("E51").Value = "  -0.46%  "
